# ---------------------------------------------------------------------------
# Issues fixes and reports
#
# AMSIN : append rows 63-64
# BETA  : append row 33
# AMS   : fix formatting on row 36 (was missing explicit styles) + append row 37
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Helper row used as a scratch "donor" cell so that date-looking strings can be
# written as literal text (bypassing Excel's automatic date recognition) while
# still landing on the same number format/style used by the rest of the column.
$HELPER_ROW = 500

function Write-DataRow {
    param(
        $ws,
        [int]$row,
        [string]$runDate,
        [double]$runTime,
        [string]$sprintName,
        [double]$total,
        [double]$pass,
        [double]$fail,
        [double]$timeTaken,
        [int]$dateFormatDonorRow,
        [int]$timeFormatDonorRow
    )

    # --- Column A: "Run Date" is stored as literal text (e.g. "2023-02-17"),
    # not as an Excel date serial. Going through a text-formatted scratch
    # cell and pasting values+formats keeps it literal text.
    $ws.Cells.Item($HELPER_ROW, 1).NumberFormat = "@"
    $ws.Cells.Item($HELPER_ROW, 1).Value = $runDate
    $ws.Cells.Item($HELPER_ROW, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4163)

    # --- Column B: "Run Time" is a real number (date+time serial). Copy the
    # number format from an existing cell in the column first, then assign
    # the value directly so the existing style id is reused.
    $ws.Cells.Item($dateFormatDonorRow, 2).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
    $ws.Cells.Item($row, 2).Value = $runTime

    # --- Column C: Sprint name, plain text (not date-like) - safe to assign directly.
    $ws.Cells.Item($row, 3).Value = $sprintName

    # --- Columns D-G: plain numbers.
    $ws.Cells.Item($row, 4).Value = $total
    $ws.Cells.Item($row, 5).Value = $pass
    $ws.Cells.Item($row, 6).Value = $fail
    $ws.Cells.Item($row, 7).Value = $timeTaken

    $ws.Cells.Item($HELPER_ROW, 1).Clear()
}

# ---------------------------------------------------------------------------
# AMSIN - add rows 63 and 64
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Write-DataRow -ws $wsAmsin -row 63 -runDate "2023-02-17" -runTime 44974.44220613426 `
    -sprintName "173cyclefst" -total 75 -pass 73 -fail 2 -timeTaken 2.11 `
    -dateFormatDonorRow 62 -timeFormatDonorRow 62

Write-DataRow -ws $wsAmsin -row 64 -runDate "2023-02-20" -runTime 44977.40638398148 `
    -sprintName "173fnlrun" -total 75 -pass 74 -fail 1 -timeTaken 1.82 `
    -dateFormatDonorRow 63 -timeFormatDonorRow 63

# ---------------------------------------------------------------------------
# BETA - add row 33
# ---------------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

Write-DataRow -ws $wsBeta -row 33 -runDate "2023-02-20" -runTime 44977.58886039352 `
    -sprintName "173beta" -total 75 -pass 75 -fail 0 -timeTaken 2.26 `
    -dateFormatDonorRow 32 -timeFormatDonorRow 32

# ---------------------------------------------------------------------------
# AMS - append row 37 FIRST (it must copy row 36's current/pre-fix styling,
# which at this point still has no explicit style applied - matching the
# target file where the newly appended row keeps that same "unstyled" look),
# then fix up row 36's own styling + refine its run-time value.
# ---------------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# copy row 36's current (pre-fix) per-cell formatting onto row 37 before we
# touch row 36 at all.
$wsAms.Cells.Item(36, 1).Copy()
$wsAms.Cells.Item(37, 1).PasteSpecial(-4122)
$wsAms.Cells.Item(36, 2).Copy()
$wsAms.Cells.Item(37, 2).PasteSpecial(-4122)
$wsAms.Cells.Item(36, 3).Copy()
$wsAms.Cells.Item(37, 3).PasteSpecial(-4122)
$wsAms.Cells.Item(36, 4).Copy()
$wsAms.Cells.Item(37, 4).PasteSpecial(-4122)
$wsAms.Cells.Item(36, 5).Copy()
$wsAms.Cells.Item(37, 5).PasteSpecial(-4122)
$wsAms.Cells.Item(36, 6).Copy()
$wsAms.Cells.Item(37, 6).PasteSpecial(-4122)
$wsAms.Cells.Item(36, 7).Copy()
$wsAms.Cells.Item(37, 7).PasteSpecial(-4122)

# column A (date-as-text) needs the text-scratch-cell trick so it isn't
# reinterpreted as a date; re-apply row 36's (unstyled) look afterwards.
$wsAms.Cells.Item($HELPER_ROW, 1).NumberFormat = "@"
$wsAms.Cells.Item($HELPER_ROW, 1).Value = "2023-02-20"
$wsAms.Cells.Item($HELPER_ROW, 1).Copy()
$wsAms.Cells.Item(37, 1).PasteSpecial(-4163)
$wsAms.Cells.Item(36, 1).Copy()
$wsAms.Cells.Item(37, 1).PasteSpecial(-4122)
$wsAms.Cells.Item($HELPER_ROW, 1).Clear()

$wsAms.Cells.Item(37, 2).Value = 44977.83604617552
$wsAms.Cells.Item(37, 3).Value = "live173"
$wsAms.Cells.Item(37, 4).Value = 75
$wsAms.Cells.Item(37, 5).Value = 75
$wsAms.Cells.Item(37, 6).Value = 0
$wsAms.Cells.Item(37, 7).Value = 2.7

# Now fix row 36's own styling: give it the same explicit "General" style the
# rest of the sheet uses for columns A, C-G (it previously had none), and
# refine the stored run-time value's precision.
$wsAms.Cells.Item(36, 1).NumberFormat = "General"
$wsAms.Cells.Item(36, 3).NumberFormat = "General"
$wsAms.Cells.Item(36, 4).NumberFormat = "General"
$wsAms.Cells.Item(36, 5).NumberFormat = "General"
$wsAms.Cells.Item(36, 6).NumberFormat = "General"
$wsAms.Cells.Item(36, 7).NumberFormat = "General"
$wsAms.Cells.Item(36, 2).Value = 44946.89340409722
